$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle existing row 9 (it is no longer the start of the next block;
#     it becomes a block-closing row like rows 4 and 7) ---
$ws.Range("A7:E7").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)

# --- New row 10 (new one-row block, style like row 5) ---
$ws.Range("A5:E5").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Cells.Item(10,1).Value = "SCRIPT/G01P03A/um1413.ssb"
$ws.Cells.Item(10,2).Value = 160
$ws.Cells.Item(10,3).Value = " There\'s no way we can ignore\nwhat [CS:N]Grovyle[CR] is up to."
$ws.Cells.Item(10,4).Value = " Мы не можем игнорировать то,\nчто затевает [CS:N]Гровайл[CR]."
$ws.Cells.Item(10,5).Value = " Íú îå íïçåí éãîïñéñïâàóû óï,\nœóï èàóåâàåó [CS:N]Ãñïâàêì[CR]."

# --- New row 11 (first row of a new two-row block, style like row 6) ---
$ws.Range("A6:E6").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Cells.Item(11,1).Value = "SCRIPT/G01P03A/um1606.ssb"
$ws.Cells.Item(11,2).Value = 138
$ws.Cells.Item(11,3).Value = " But if [CS:N]Grovyle[CR] is out to get the\nTime Gears...[K]it might be better to just seal\nthem away."
$ws.Cells.Item(11,4).Value = " Но если [CS:N]Гровайл[CR] будет пытаться\nукрасть Шестерни Времени...[K] Возможно, нам\nстоит их спрятать."
$ws.Cells.Item(11,5).Value = " Îï åòìé [CS:N]Ãñïâàêì[CR] áôäåó ðúóàóûòÿ\nôëñàòóû Šåòóåñîé Âñåíåîé...[K] Âïèíïçîï, îàí\nòóïéó éö òðñÿóàóû."

# --- New row 12 (continuation row of the row-11 block, style like row 3/8, no A cell) ---
$ws.Range("B3:E3").Copy()
$ws.Range("B12:E12").PasteSpecial(-4122)
$ws.Rows.Item(12).RowHeight = 21.6
$ws.Cells.Item(12,2).Value = 141
$ws.Cells.Item(12,3).Value = " It\'s an option worth considering."
$ws.Cells.Item(12,4).Value = " Над таким вариантом стоит\nпоразмыслить."
$ws.Cells.Item(12,5).Value = " Îàä óàëéí âàñéàîóïí òóïéó\nðïñàèíúòìéóû."

$excel.CutCopyMode = $false

# --- Update worksheet dimension / view to match ---
$ws.Range("A1:E12").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D11").Select()

